$wb = $excel.ActiveWorkbook

# ===== Sheet 展览 (Exhibitions): F column 'want to go' count bumps =====
$ws1 = $wb.Worksheets.Item('展览')
$ws1.Range("F2").Value = 2455
$ws1.Range("F5").Value = 411
$ws1.Range("F6").Value = 686
$ws1.Range("F9").Value = 565
$ws1.Range("F10").Value = 929
$ws1.Range("F12").Value = 124
$ws1.Range("F14").Value = 50
$ws1.Range("F16").Value = 1068
$ws1.Range("F17").Value = 23933
$ws1.Range("F18").Value = 2204
$ws1.Range("F19").Value = 141
$ws1.Range("F22").Value = 47
$ws1.Range("F23").Value = 349
$ws1.Range("F25").Value = 64
$ws1.Range("F28").Value = 52
$ws1.Range("F30").Value = 342

# ===== Sheet 演出 (Performances): F column bumps =====
$ws2 = $wb.Worksheets.Item('演出')
$ws2.Range("F4").Value = 189
$ws2.Range("F7").Value = 253
$ws2.Range("F8").Value = 20
$ws2.Range("F10").Value = 3596
$ws2.Range("F12").Value = 145
$ws2.Range("F16").Value = 17
$ws2.Range("F17").Value = 134
$ws2.Range("F19").Value = 4114

# ===== Sheet 本地生活 (Local life): F column bumps =====
$ws3 = $wb.Worksheets.Item('本地生活')
$ws3.Range("F3").Value = 160
$ws3.Range("F4").Value = 747
$ws3.Range("F5").Value = 239

# ===== Sheet 全部类型 (All types): F column bumps (rows whose event identity is unchanged) =====
$ws4 = $wb.Worksheets.Item('全部类型')
$ws4.Range("F3").Value = 160
$ws4.Range("F4").Value = 2455
$ws4.Range("F5").Value = 747
$ws4.Range("F8").Value = 411
$ws4.Range("F9").Value = 686
$ws4.Range("F11").Value = 189
$ws4.Range("F14").Value = 253
$ws4.Range("F15").Value = 239
$ws4.Range("F17").Value = 565
$ws4.Range("F18").Value = 929
$ws4.Range("F19").Value = 124
$ws4.Range("F21").Value = 50
$ws4.Range("F23").Value = 1068
$ws4.Range("F24").Value = 23933
$ws4.Range("F25").Value = 20
$ws4.Range("F28").Value = 145
$ws4.Range("F30").Value = 2204
$ws4.Range("F31").Value = 141
$ws4.Range("F36").Value = 349

# ===== Sheet 全部类型 rows 38-48: event list update =====
# The 2024-08-04 'wio夏时之鸢代号鸢Only' event (old row 38) was removed from the
# feed; events previously on rows 39-48 shift up to rows 38-47, and a brand-new event
# (LoveLive! 10th-anniversary tour) is appended as the new row 48. Column A (the static
# index) is left untouched throughout, matching upstream's diff (only B..I move).
#
# Column B holds plain-text dates like '2024-08-04'; a direct .Value assignment of such
# a string gets auto-coerced to a real date serial by Excel's smart input (and stamps a
# NumberFormat style on the cell). To preserve the original plain-text storage exactly,
# we stage the literal text as a formula result in a scratch cell and paste-special just
# the value over the target -- paste-special values do not re-run text->date coercion.
$scratch = $ws4.Range("K200")
$scratch.Formula = "='2024-08-04'"
$scratch.Copy()
$ws4.Range("B38").PasteSpecial(-4163)
$ws4.Range("C38").Value = '广州·星之光动漫嘉年华'
$ws4.Range("D38").Value = '钟村镇105国道西侧 广州雄峰城展览中心'
$ws4.Range("E38").Value = '2024.08.04 10:00-08.04 17:00'
$ws4.Range("F38").Value = 226
$ws4.Range("G38").Value = 60
$ws4.Range("H38").Value = 'https://show.bilibili.com/platform/detail.html?id=87077'
$ws4.Range("I38").Value = '//i2.hdslb.com/bfs/openplatform/202406/hOZ6VVFx1717571239392.jpeg'

$scratch.Formula = "='2024-08-04'"
$scratch.Copy()
$ws4.Range("B39").PasteSpecial(-4163)
$ws4.Range("C39").Value = '广州·格斗游戏FTGonly'
$ws4.Range("D39").Value = '芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋'
$ws4.Range("E39").Value = '2024.08.04 10:00-08.04 19:00'
$ws4.Range("F39").Value = 15
$ws4.Range("G39").Value = 68
$ws4.Range("H39").Value = 'https://show.bilibili.com/platform/detail.html?id=87090'
$ws4.Range("I39").Value = '//i1.hdslb.com/bfs/openplatform/202406/Vk8sR8Oj1717582522018.png'

$scratch.Formula = "='2024-08-10'"
$scratch.Copy()
$ws4.Range("B40").PasteSpecial(-4163)
$ws4.Range("C40").Value = '广州·系统任务：重生之我是音乐一体机！王子健2024巡回演出'
$ws4.Range("D40").Value = '广州天河区花城大道89号美食街北二门 SD Livehouse'
$ws4.Range("E40").Value = '2024.08.10 20:00-08.10 22:00'
$ws4.Range("F40").Value = 15
$ws4.Range("G40").Value = 328
$ws4.Range("H40").Value = 'https://show.bilibili.com/platform/detail.html?id=87585'
$ws4.Range("I40").Value = '//i0.hdslb.com/bfs/openplatform/202406/zIb7ZnHb1718675848837.jpeg'

$scratch.Formula = "='2024-08-11'"
$scratch.Copy()
$ws4.Range("B41").PasteSpecial(-4163)
$ws4.Range("C41").Value = '广州·咒术回战ONLY'
$ws4.Range("D41").Value = '西环路1号 广州岭南会展中心'
$ws4.Range("E41").Value = '2024.08.11 10:00-08.11 17:00'
$ws4.Range("F41").Value = 49
$ws4.Range("G41").Value = 60
$ws4.Range("H41").Value = 'https://show.bilibili.com/platform/detail.html?id=87433'
$ws4.Range("I41").Value = '//i1.hdslb.com/bfs/openplatform/202406/kNv9yqGn1718350051848.jpeg'

$scratch.Formula = "='2024-08-11'"
$scratch.Copy()
$ws4.Range("B42").PasteSpecial(-4163)
$ws4.Range("C42").Value = '广州·妖都原神&崩铁only-清凉大作战-'
$ws4.Range("D42").Value = '黄边三横路一街1号 设计殿堂'
$ws4.Range("E42").Value = '2024.08.11 10:00-08.11 16:30'
$ws4.Range("F42").Value = 44
$ws4.Range("G42").Value = 60
$ws4.Range("H42").Value = 'https://show.bilibili.com/platform/detail.html?id=87321'
$ws4.Range("I42").Value = '//i1.hdslb.com/bfs/openplatform/202406/7k54Bi4X1718025336899.jpeg'

$scratch.Formula = "='2024-08-14'"
$scratch.Copy()
$ws4.Range("B43").PasteSpecial(-4163)
$ws4.Range("C43").Value = '广州·Marcin Patrzalek 2024 《原声之龙》指弹吉他音乐会'
$ws4.Range("D43").Value = '海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse'
$ws4.Range("E43").Value = '2024.08.14 20:00-08.14 21:30'
$ws4.Range("F43").Value = 133
$ws4.Range("G43").Value = 380
$ws4.Range("H43").Value = 'https://show.bilibili.com/platform/detail.html?id=86291'
$ws4.Range("I43").Value = '//i1.hdslb.com/bfs/openplatform/202405/vsOXym1L1716546835148.jpeg'

$scratch.Formula = "='2024-08-16'"
$scratch.Copy()
$ws4.Range("B44").PasteSpecial(-4163)
$ws4.Range("C44").Value = '广州·《最后的莫西干人》-印第安音乐家亚历桑德罗&丛林回响乐队巡演'
$ws4.Range("D44").Value = '东风中路299号 广州中山纪念堂'
$ws4.Range("E44").Value = '2024.08.16 20:00-08.16 21:30'
$ws4.Range("F44").Value = 2
$ws4.Range("G44").Value = 380
$ws4.Range("H44").Value = 'https://show.bilibili.com/platform/detail.html?id=86143'
$ws4.Range("I44").Value = '//i1.hdslb.com/bfs/openplatform/202405/4oOXA1j01716175554059.jpeg'

$scratch.Formula = "='2024-08-17'"
$scratch.Copy()
$ws4.Range("B45").PasteSpecial(-4163)
$ws4.Range("C45").Value = '广州·鸟山明作品《龙珠》40周年only纪念展'
$ws4.Range("D45").Value = '逸景路462号珠江国际纺织城d区6层 珠江时尚馆'
$ws4.Range("E45").Value = '2024.08.17 10:00-08.17 17:30'
$ws4.Range("F45").Value = 25
$ws4.Range("G45").Value = 68
$ws4.Range("H45").Value = 'https://show.bilibili.com/platform/detail.html?id=86780'
$ws4.Range("I45").Value = '//i1.hdslb.com/bfs/openplatform/202405/4k7Thger1717147185584.jpeg'

$scratch.Formula = "='2024-08-18'"
$scratch.Copy()
$ws4.Range("B46").PasteSpecial(-4163)
$ws4.Range("C46").Value = '广州·原神×崩坏×绝区零only'
$ws4.Range("D46").Value = '西环路1号 广州岭南会展中心'
$ws4.Range("E46").Value = '2024.08.18 10:00-08.18 17:00'
$ws4.Range("F46").Value = 432
$ws4.Range("G46").Value = 60
$ws4.Range("H46").Value = 'https://show.bilibili.com/platform/detail.html?id=87025'
$ws4.Range("I46").Value = '//i0.hdslb.com/bfs/openplatform/202405/lsOq4H701717169339283.png'

$scratch.Formula = "='2024-08-23'"
$scratch.Copy()
$ws4.Range("B47").PasteSpecial(-4163)
$ws4.Range("C47").Value = '广州·LoveLiveOnly'
$ws4.Range("D47").Value = '芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋'
$ws4.Range("E47").Value = '2024.08.23 10:00-08.23 19:00'
$ws4.Range("F47").Value = 186
$ws4.Range("G47").Value = 68.8
$ws4.Range("H47").Value = 'https://show.bilibili.com/platform/detail.html?id=87033'
$ws4.Range("I47").Value = '//i2.hdslb.com/bfs/openplatform/202406/a8shiH411717579829497.jpeg'

$scratch.Formula = "='2024-08-24'"
$scratch.Copy()
$ws4.Range("B48").PasteSpecial(-4163)
$ws4.Range("C48").Value = '广州·LoveLive！电视动画播放十周年纪念巡演'
$ws4.Range("D48").Value = '机场路1733号 久米空间LIVEHOUSE'
$ws4.Range("E48").Value = '2024.08.24 12:30-08.25 18:30'
$ws4.Range("F48").Value = 4114
$ws4.Range("G48").Value = '已售罄'
$ws4.Range("H48").Value = 'https://show.bilibili.com/platform/detail.html?id=86959'
$ws4.Range("I48").Value = '//i1.hdslb.com/bfs/openplatform/202406/apzqBc5d1717661406596.jpeg'

$scratch.Clear()

